# Update column A values for several rows (algorithm name / data update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.879
$ws.Range("A14").Value = -21.66
$ws.Range("A16").Value = -22.163
$ws.Range("A21").Value = -20.022
$ws.Range("A23").Value = -20.203
$ws.Range("A25").Value = -21.632
$ws.Range("A26").Value = -21.692
$ws.Range("A29").Value = -21.086
$ws.Range("A40").Value = -20
$ws.Range("A53").Value = -21.853
$ws.Range("A57").Value = -22.219
$ws.Range("A59").Value = -22.607
$ws.Range("A65").Value = -21.529
$ws.Range("A69").Value = -21.519
$ws.Range("A79").Value = -21.087
$ws.Range("A83").Value = -21.938
$ws.Range("A91").Value = -20.666
$ws.Range("A93").Value = -21.508
$ws.Range("A100").Value = -22.277
$ws.Range("A103").Value = -22.086
